# "Updated examples / scripts till 24-08-2025"
#
# Changes applied (content-level, i.e. the things actually reachable through
# the Excel object model -- file-format bookkeeping such as fileVersion/
# rupBuild, xr:revisionPtr, the x15ac:absPath session URL, the bookViews
# window pixel size/xr2:uid, and the metadata.xml xlrd namespace are all
# Excel-session/runtime artifacts stamped on save, not user-settable state,
# so they are intentionally left alone here):
#
#   1. Rename the second sheet:
#        "SCAN accmu Real UseCase" -> "SCAN To workingout balance"
#   2. On sheet 1 ("SCAN To accumulate With Reset"): it becomes the
#      selected/active tab, with B2:B32 selected (active cell B2).
#   3. On sheet 2 ("SCAN To workingout balance"): it is no longer the
#      active tab, and its last selection is E23.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Rename sheet 2
$ws2.Name = "SCAN To workingout balance"

# 2 & 3. Set each sheet's remembered selection, then leave sheet 1 active
#    (order matters: the last .Select()'d sheet/range becomes the
#    tabSelected / workbook's active tab).
$ws2.Select()
$ws2.Range("E23").Select()

$ws1.Select()
$ws1.Range("B2:B32").Select()
